# Update dashboards - 2025-12-13
# Applies the latest data pull to the rate/date cells in the dashboard table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Aguilar Prototype")

# Row 29 - T5YIFR (5yr, 5yr Forward)
$ws.Range("N29").Value = 46003
$ws.Range("Q29").Value = 2.2
$ws.Range("S29").Value = 2.18
$ws.Range("U29").Value = 2.2

# Row 30 - T10YIE (10yr TIPS)
$ws.Range("N30").Value = 46003
$ws.Range("Q30").Value = 2.26
$ws.Range("S30").Value = 2.25
$ws.Range("U30").Value = 2.26

# Row 47 - DFF (FFR)
$ws.Range("N47").Value = 46002
$ws.Range("Q47").Value = 3.64

# Row 48 - DGS2 (2y UST)
$ws.Range("N48").Value = 46002
$ws.Range("Q48").Value = 3.52
$ws.Range("R48").Value = 3.54
$ws.Range("S48").Value = 3.61
$ws.Range("T48").Value = 3.57

# Row 49 - DGS5 (5y UST)
$ws.Range("N49").Value = 46002
$ws.Range("R49").Value = 3.72
$ws.Range("S49").Value = 3.78
$ws.Range("T49").Value = 3.75

# Row 50 - DGS10 (10y UST)
$ws.Range("N50").Value = 46002
$ws.Range("Q50").Value = 4.14
$ws.Range("R50").Value = 4.13
$ws.Range("S50").Value = 4.18
$ws.Range("T50").Value = 4.17

# Row 52 - DBAA (BAA)
$ws.Range("N52").Value = 46002
$ws.Range("Q52").Value = 5.87
$ws.Range("R52").Value = 5.89
$ws.Range("S52").Value = 5.91
$ws.Range("T52").Value = 5.9
